$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 9
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 7
$ws.Range("B4").Value = 1.1000000000000001
$ws.Range("C5").Value = 16

# Update selection to match new active cell
$ws.Range("D9").Select()

# Update window size (matches the new windowWidth/windowHeight recorded
# in the workbook view when the file was last saved)
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 12180

